# Remove column from alcohol measurement data.
#
# The "M" column (13) of Sheet1 held duplicate/old data; it is removed and
# the following column ("N") shifts left to take its place, so the sheet
# ends up with columns A:M instead of A:N.

$wb = $excel.ActiveWorkbook

# Sheet1 holds the actual alcohol measurement data - delete column M so that
# the old column N values become the new column M values.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Columns.Item(13).Delete()

# Reset the zoom level on every sheet (matches the resaved workbook).
foreach ($ws in $wb.Worksheets) {
    $ws.Activate()
    $excel.ActiveWindow.Zoom = 95
}

# Restore Sheet1 as the active/selected sheet, with the selection resting
# on the new last column (M1), mirroring where the deleted column used to be.
$ws1.Activate()
$ws1.Range("M1").Select() | Out-Null
